$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$b = $ws.Range("D1").Borders.Item(9)
$b.LineStyle = 1
$b.ColorIndex = -4105
Write-Host "done"
